$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit -----------------------------------------------------------
# F1 used to hold the formula "=B7" (which evaluated to 59.8). The author
# replaced it with a plain literal value of 60.1, which ripples through the
# dependent formulas in F3:F21 (signal model), H3:H21 (inverse power calc),
# J3:J21 (percentage) and K3:K21 / K23 (error metrics), as well as the
# "Computed RSSI" series cached in chart 4.
$ws.Range("F1").Value = 60.1

# --- Chart reposition ------------------------------------------------------
# The 4th chart ("Gr\u00e1fico 4") was dragged to a new spot on the sheet
# (up/left of its previous location), while keeping essentially the same
# size. Values below are the point-based Left/Top/Width/Height that produce
# the target anchor (from col 12 / row 3 to col 23 / row 24).
$chartObj = $ws.ChartObjects().Item(4)
$chartObj.Left = 773.2162109375
$chartObj.Top = 48
$chartObj.Width = 626.9125
$chartObj.Height = 322.8

# --- Selection --------------------------------------------------------------
# The author finished up with I3:K21 selected (active cell I3).
$ws.Range("I3:K21").Select() | Out-Null
